# Fill in new execution metrics for configuration 11 (UC3.2_TC1)
# on the "QuantitativeMetrics" sheet, "Execution metrics" table (rows 5-7):
#   Compilation success    -> Value: yes,  Note: (cleared)
#   Runtime without error  -> Value: yes
#   Assertion validity     -> Value: no,   Note: invalid expect
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("QuantitativeMetrics")

$ws.Range("B5").Value = "yes"
$ws.Range("C5").Value = ""

$ws.Range("B6").Value = "yes"

$ws.Range("B7").Value = "no"
$ws.Range("C7").Value = "invalid expect"
